$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Small text tweaks to existing rows
# ---------------------------------------------------------------------------
# C8 gained a "Uses" description that was previously blank.
$ws.Range("C8").Value = "To show remote origin URL"

# B9 command now wraps the URL in single quotes.
$ws.Range("B9").Value = "git remote set-url origin 'https:// PAT TOKEN @git repository address'"

# ---------------------------------------------------------------------------
# 2) Finish row 14 (previously only had the S.No. filled in) and append the
#    23 new command rows (15-24), extending the table down to row 24.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "git restore --staged file_name"
$ws.Range("C14").Value = "for rollback from staged to untracked area"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "git restore file name"
$ws.Range("C15").Value = "To restore deleted file from git "

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "rm file name"
$ws.Range("C16").Value = "for remove the file "

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "git checkout any branch"
$ws.Range("C17").Value = "To switch on perticular branch"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "git checkout -b Any new branch name "
$ws.Range("C18").Value = "To Create a new from from this existing branch and switch on new branch"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "git branch -d any branch name"
$ws.Range("C19").Value = "To remove any perticular branch"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "git fetch"
$ws.Range("C20").Value = "To fetch the all branches from remote to local"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "git remote add origin (remote git URL)"
$ws.Range("C21").Value = "To add remote Origin Url"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "git remote remove origin"
$ws.Range("C22").Value = "To remove remote origin url"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "git clone "

$ws.Range("A24").Value = 23
$sshText = 'GIT_SSH_COMMAND="ssh -i /home/ubuntu/.ssh/ your private key" git clone your ssh clone url address'
$ws.Range("B24").Value = $sshText
$ws.Range("B24").Characters(43, 18).Font.Bold = $true
$ws.Range("B24").Characters(72, 26).Font.Bold = $true
$ws.Range("C24").Value = "To use any perticular ssh private key instead of Default ssh private key for ssh git clone"

# ---------------------------------------------------------------------------
# 3) Auto-fit columns B & C now that the longer command text has been added.
# ---------------------------------------------------------------------------
$ws.Columns("B:C").AutoFit()

# ---------------------------------------------------------------------------
# 4) Threaded comment on C15 describing the git status / git restore flow.
# ---------------------------------------------------------------------------
$ws.Range("C15").AddCommentThreaded("First check our deleted file are showing are not we can check it with (git status) command if we find deleted file is showing here then we can restore it with git restore command")

# ---------------------------------------------------------------------------
# 5) Misc view tweaks that show up in the diff.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("C10").Select()
